# Reorder the header row (row 1) of the MimsSoil sheet, moving the
# samp_taxon_id / rel_to_oxygen / ... headers (and their associated
# dataValidation list ranges) to new column positions, per the
# "more write_mixs_linkml param updates" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MimsSoil")

# New header order for columns A..AR (columns AS onward are unchanged).
$headers = @(
    "ref_biomaterial",
    "samp_collect_method",
    "sim_search_meth",
    "project_name",
    "elev",
    "samp_collect_device",
    "tax_class",
    "mid",
    "depth",
    "adapters",
    "neg_cont_type",
    "assembly_name",
    "samp_taxon_id",
    "samp_name",
    "lib_reads_seqd",
    "assembly_qual",
    "geo_loc_name",
    "annot",
    "collection_date",
    "size_frac",
    "lib_layout",
    "nucl_acid_ext",
    "rel_to_oxygen",
    "lat_lon",
    "env_local_scale",
    "samp_vol_we_dna_ext",
    "assembly_software",
    "samp_size",
    "temp",
    "samp_mat_process",
    "ref_db",
    "nucl_acid_amp",
    "feat_pred",
    "env_broad_scale",
    "lib_screen",
    "env_medium",
    "number_contig",
    "pos_cont_type",
    "lib_vector",
    "source_mat_id",
    "experimental_factor",
    "seq_meth",
    "alt",
    "lib_size"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# The three dataValidation list ranges that track these moved headers need
# to move along with them: rel_to_oxygen (B->W), neg_cont_type (M->K),
# lib_layout (AE->U). Rebuild all 8 validations (3 moved + 5 unchanged) in
# the new order so the sheet's <dataValidations> sequence matches the
# target file exactly.
function Add-ListValidation($range, $listFormula) {
    $v = $ws.Range($range).Validation
    $v.Add(3, 1, 1, $listFormula)
    $v.InCellDropdown = $true
    $v.ShowInput = $false
    $v.ShowError = $false
}

$ws.Cells.Validation.Delete()

Add-ListValidation "K2:K1048576" '"DNA-free PCR mix,distilled water,empty collection device,empty collection tube,phosphate buffer,sterile swab,sterile syringe"'
Add-ListValidation "U2:U1048576" '"other,paired,single,vector"'
Add-ListValidation "W2:W1048576" '"aerobe,anaerobe,facultative,microaerophilic,microanaerobe,obligate aerobe,obligate anaerobe"'
Add-ListValidation "BB2:BB1048576" '"chisel,cutting disc,disc plough,drill,mouldboard,ridge till,strip tillage,tined,zonal tillage"'
Add-ListValidation "BF2:BF1048576" '"A horizon,B horizon,C horizon,E horizon,O horizon,Permafrost,R layer"'
Add-ListValidation "BS2:BS1048576" '"Acrisols,Andosols,Arenosols,Cambisols,Chernozems,Ferralsols,Fluvisols,Gleysols,Greyzems,Gypsisols,Histosols,Kastanozems,Lithosols,Luvisols,Nitosols,Phaeozems,Planosols,Podzols,Podzoluvisols,Rankers,Regosols,Rendzinas,Solonchaks,Solonetz,Vertisols,Yermosols"'
Add-ListValidation "CA2:CA1048576" '"backslope,footslope,shoulder,summit,toeslope"'
Add-ListValidation "CB2:CB1048576" '"excessively drained,moderately well,poorly,somewhat poorly,very poorly,well"'

$wb.Save()
